# Updates the "Estado de Cuenta" worksheet with the refreshed account-statement
# database: replaces the debtor/period detail rows (16-26) with the new
# part-1 data set (3 workers x periods 2003/2004), refreshes the summary
# totals (Valor Mora, Cant. Periodos) and removes the now-unused trailing
# rows so the signature block moves up directly beneath the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary header figures -------------------------------------------------
# Valor Mora (total) and Cant. Periodos reflect the new, smaller data set.
$ws.Range("E11").Value2 = 595112
$ws.Range("F13").Value2 = 2

# --- Detail rows 16-20: new data (3 trabajadores x 2 periodos) -------------
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "80163620"
$ws.Range("D16").Value2 = "ROQUE PEDROZA NAVARRO"
$ws.Range("E16").Value2 = "2003"
$ws.Range("F16").Value2 = 120000
$ws.Range("G16").Value2 = 877803

$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "1047416123"
$ws.Range("D17").Value2 = "SUSANA MARGARITA CANCHILA DE LA ESPRIELLA"
$ws.Range("E17").Value2 = "2003"
$ws.Range("F17").Value2 = 160000
$ws.Range("G17").Value2 = 877803

$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "1048941338"
$ws.Range("D18").Value2 = "DANIEL EDUARDO DORIA MELENDEZ"
$ws.Range("E18").Value2 = "2003"
$ws.Range("F18").Value2 = 35112
$ws.Range("G18").Value2 = 877803

$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "80163620"
$ws.Range("D19").Value2 = "ROQUE PEDROZA NAVARRO"
$ws.Range("E19").Value2 = "2004"
$ws.Range("F19").Value2 = 120000
$ws.Range("G19").Value2 = 877803

$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "1047416123"
$ws.Range("D20").Value2 = "SUSANA MARGARITA CANCHILA DE LA ESPRIELLA"
$ws.Range("E20").Value2 = "2004"
$ws.Range("F20").Value2 = 160000
$ws.Range("G20").Value2 = 877803

# Row 20 is now the last row of the table, so it needs the table's "bottom
# border" formatting that used to live on row 26 (the old last row). Copy
# that formatting across before the old rows are removed.
$ws.Range("B26:J26").Copy() | Out-Null
$ws.Range("B20:J20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Remove the old trailing rows (21-26 held the rest of the previous
# data table) so the signature block (previously rows 31-32) shifts up to
# sit right under the new, shorter table at rows 25-26. ----------------------
$ws.Range("21:26").EntireRow.Delete() | Out-Null
